# Apply "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific bullet
# paragraphs - resume achievements / work-experience responsibilities -
# matching the target diff (bold + color 2C3E50 run splits).

$d = $word.ActiveDocument
$highlightColor = 5258796   # BGR long value equivalent to hex RGB 2C3E50

# Locate the paragraph whose text contains $mustContain but does NOT
# contain $mustNotContain (used to disambiguate the two near-duplicate
# "Achieved 87% ... 71%" bullets - one has a "reducing polling..." tail,
# the other doesn't).
function Find-TargetParagraph {
    param(
        $Doc,
        [string]$MustContain,
        [string]$MustNotContain
    )

    foreach ($p in $Doc.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Contains($MustContain)) {
            if ($MustNotContain -eq "" -or -not $t.Contains($MustNotContain)) {
                return $p
            }
        }
    }
    return $null
}

# Bold + color every listed metric substring, left-to-right, constrained
# to the paragraph's own range so sibling paragraphs are never touched.
function Highlight-Metrics {
    param(
        $Paragraph,
        [string[]]$Metrics
    )

    if ($Paragraph -eq $null) {
        return
    }

    $pStart = $Paragraph.Range.Start
    $pEnd = $Paragraph.Range.End
    $cursor = $pStart

    foreach ($metric in $Metrics) {
        $search = $d.Range($cursor, $pEnd)
        $found = $search.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $search.Font.Bold = 1
            $search.Font.Color = $highlightColor
            $cursor = $search.End
        }
    }
}

# 1) "...developed geospatial machine learning algorithms improving
#     demographic classification accuracy from 23% to 64%"
$p1 = Find-TargetParagraph $d "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed" ""
Highlight-Metrics $p1 @("23%", "64%")

# 2) "Achieved 87% prediction accuracy ... 71%, reducing polling error
#     margins from ±4.2% to ±2.1%"
$p2 = Find-TargetParagraph $d "reducing polling error margins" ""
Highlight-Metrics $p2 @("87%", "71%", "±4.2%", "±2.1%")

# 3) "Wrote RFP and analyzed bids from 1,200 vendors..."
$p3 = Find-TargetParagraph $d "Wrote RFP and analyzed bids from" ""
Highlight-Metrics $p3 @("1,200")

# 4) "...became the $400M Polling Consortium Database ... valued at $1B+"
$p4 = Find-TargetParagraph $d "Polling Consortium Database at The Analyst Institute" ""
Highlight-Metrics $p4 @("$400M", "$1B")

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and
#     organizations $4.7M"
$p5 = Find-TargetParagraph $d "Algorithm reduced mapping costs by" ""
Highlight-Metrics $p5 @("73.5%", "$4.7M")

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry
#     standard of 71%" (no "reducing..." tail - the Key Achievements one)
$p6 = Find-TargetParagraph $d "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" "reducing"
Highlight-Metrics $p6 @("87%", "71%")

Write-Output "done"
